# Applies the LOQ4233.docx edit: splits several concatenated-sentence
# paragraphs into separate lines using manual line breaks (<w:br/>),
# matching the target OOXML diff. Uses Find/Replace with the Word
# special search code ^l (manual line break) as the replacement marker,
# which Word renders as a <w:br/> element and splits the run's <w:t>.

$d = $word.ActiveDocument

# Objetivos bullet: split 3 sentences
$rng = $d.Content
$found = $rng.Find.Execute('Apresentar ao aluno o conceito de uma organização e os fundamentos de sua administração;Caracterizar as diversas áreas funcionais existentes nas organizações;Despertar o interesse dos alunos para questões de gestão', $false, $false, $false, $false, $false, $true, 1, $false, 'Apresentar ao aluno o conceito de uma organização e os fundamentos de sua administração;^lCaracterizar as diversas áreas funcionais existentes nas organizações;^lDespertar o interesse dos alunos para questões de gestão', 2)
Write-Host "Step 1 found: $found"

# Programa (PT) paragraph: split into 4 lines
$rng = $d.Content
$found = $rng.Find.Execute('1 - A Administração das organizações - definindo a administração2 - O processo administrativo: planejamento, organização, direção, controle3 – Processos de Gestão: Marketing, Finanças, Gestão de Pessoas, Produção e Operações, Pesquisa e Desenvolvimento, Tecnologia da Informação, Logística e Meio Ambiente.A disciplina será ministrada com duas estratégias pedagógicas a) aplicação de diferentes métodos ativos para compreender os principais conceitos necessários à gestão de negócios, e b) aplicação de conceitos por meio do Programa de Aprendizagem com Extensão, por meio do qual o alunos oferecem consultoria a micro e pequenas empresas da região de Lorena ou de parentes e amigos. Nestas consultorias times de alunos, orientados pelo professor, se debruçam sobre um pequeno problema de gestão da empresa e oferecem soluções.', $false, $false, $false, $false, $false, $true, 1, $false, '1 - A Administração das organizações - definindo a administração^l2 - O processo administrativo: planejamento, organização, direção, controle^l3 – Processos de Gestão: Marketing, Finanças, Gestão de Pessoas, Produção e Operações, Pesquisa e Desenvolvimento, Tecnologia da Informação, Logística e Meio Ambiente.^lA disciplina será ministrada com duas estratégias pedagógicas a) aplicação de diferentes métodos ativos para compreender os principais conceitos necessários à gestão de negócios, e b) aplicação de conceitos por meio do Programa de Aprendizagem com Extensão, por meio do qual o alunos oferecem consultoria a micro e pequenas empresas da região de Lorena ou de parentes e amigos. Nestas consultorias times de alunos, orientados pelo professor, se debruçam sobre um pequeno problema de gestão da empresa e oferecem soluções.', 2)
Write-Host "Step 2 found: $found"

# Bibliografia paragraph: split into 3 entries (blank line between)
$rng = $d.Content
$found = $rng.Find.Execute('LEMOS, Paulo de Mattos et al. Gestão estratégica de empresas. Rio de Janeiro: Fundação Getúlio Vargas, 2014.Ludovico, Nelson. Gestão estratégica de negócios. São Paulo: Saraiva, 2018Serra, Fernando Ribeiro et al. Gestão estratégica: conceitos e casos. São Paulo: Atlas, 2014.', $false, $false, $false, $false, $false, $true, 1, $false, 'LEMOS, Paulo de Mattos et al. Gestão estratégica de empresas. Rio de Janeiro: Fundação Getúlio Vargas, 2014.^l^lLudovico, Nelson. Gestão estratégica de negócios. São Paulo: Saraiva, 2018^l^lSerra, Fernando Ribeiro et al. Gestão estratégica: conceitos e casos. São Paulo: Atlas, 2014.', 2)
Write-Host "Step 3 found: $found"

# Programa (EN) paragraph: insert break before '2 - The administrative...'
$rng = $d.Content
$found = $rng.Find.Execute(' the administration 2 - The administrati', $false, $false, $false, $false, $false, $true, 1, $false, ' the administration ^l2 - The administrati', 2)
Write-Host "Step 4 found: $found"

# Programa (EN) paragraph: insert break before '3 - Management Processes...'
$rng = $d.Content
$found = $rng.Find.Execute(' direction, control 3 - Management Proce', $false, $false, $false, $false, $false, $true, 1, $false, ' direction, control ^l3 - Management Proce', 2)
Write-Host "Step 5 found: $found"

# Programa (EN) paragraph: insert break before 'The course will be taught...'
$rng = $d.Content
$found = $rng.Find.Execute('ics and Environment. The course will be ', $false, $false, $false, $false, $false, $true, 1, $false, 'ics and Environment.^l The course will be ', 2)
Write-Host "Step 6 found: $found"

# Avaliacao / Criterio bullet: split into 2 lines
$rng = $d.Content
$found = $rng.Find.Execute('- Contribuir para a gestão de organizações de pequeno e médio porte visando melhoria de rendas de comundades;- contribuir para capacitar gestores de organizações de pequeno e medio porte.', $false, $false, $false, $false, $false, $true, 1, $false, '- Contribuir para a gestão de organizações de pequeno e médio porte visando melhoria de rendas de comundades;^l- contribuir para capacitar gestores de organizações de pequeno e medio porte.', 2)
Write-Host "Step 7 found: $found"

# Bibliografia-labeled paragraph: split into 7 lines
$rng = $d.Content
$found = $rng.Find.Execute('- Estabelecimento da comunicação aberta entre estudantes, grupo social e professor;- Acompanhamento pelo professor e grupo social da atividade a ser desenvolvida pelos alunos;- Exposição de cada grupo, sobre a proposta, desenvolvimento e finalização do projeto;- Realização de avaliação conjunta dos resultados alcançados durante a atividade, incluindo benefícios   obtidos, lições aprendidas e desafios enfrentados;- Conduzir sessões de discussão para revisar os resultados e identificar oportunidades de aplicação  prática;- Apresentação do projeto final desenvolvido para grupo social;- Avaliação do projeto apresentado, pelo grupo social e professor.', $false, $false, $false, $false, $false, $true, 1, $false, '- Estabelecimento da comunicação aberta entre estudantes, grupo social e professor;^l- Acompanhamento pelo professor e grupo social da atividade a ser desenvolvida pelos alunos;^l- Exposição de cada grupo, sobre a proposta, desenvolvimento e finalização do projeto;^l- Realização de avaliação conjunta dos resultados alcançados durante a atividade, incluindo benefícios   obtidos, lições aprendidas e desafios enfrentados;^l- Conduzir sessões de discussão para revisar os resultados e identificar oportunidades de aplicação  prática;^l- Apresentação do projeto final desenvolvido para grupo social;^l- Avaliação do projeto apresentado, pelo grupo social e professor.', 2)
Write-Host "Step 8 found: $found"
